# The "Saldo" (balance) for account 004119016 (HEMAT) was corrected from
# 62548.23 to 548.23. Because the sheet's data rows (A2:C307) are kept
# sorted by Saldo descending, fixing the value moves that account's row
# further down the list; re-sorting on column C reproduces that reordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 3).Value = 548.23

$sortRange = $ws.Range("A2:C307")
$sortKey = $ws.Range("C2:C307")
$sortRange.Sort($sortKey, 2)
